$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update headers: A1 "Fecha", B1 "Ventas", add C1 "Kpi"
$ws.Range("A1").Value = "Fecha"
$ws.Range("B1").Value = "Ventas"
$ws.Range("C1").Value = "Kpi"

# Add formula in C2 and format it as percentage
$ws.Range("C2").Formula = "=B2+B3/3"
$ws.Range("C2").NumberFormat = "0.00%"

# Select C3 to match the saved selection state
$ws.Range("C3").Select()
